$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 1830.9
$ws.Range("I32").Value = 1200.25
$ws.Range("K32").Value = 1200.25
$ws.Range("M32").Value = -874.25
# Row 98
$ws.Range("H98").Value = 1016656.06
$ws.Range("I98").Value = 1016656.06
$ws.Range("K98").Value = 1016656.06
$ws.Range("M98").Value = -1015158.06
# Row 122
$ws.Range("H122").Value = 1016656.06
$ws.Range("I122").Value = 1016656.06
$ws.Range("K122").Value = 3049968.18
$ws.Range("M122").Value = -3047518.18
# Row 132
$ws.Range("H132").Value = 352172.84
$ws.Range("I132").Value = 529018.5600000001
$ws.Range("K132").Value = 1587055.68
$ws.Range("M132").Value = -1584525.68
# Row 133
$ws.Range("H133").Value = 12200
$ws.Range("J133").Value = 12200
$ws.Range("L133").Value = 12200
$ws.Range("N133").Value = -22320

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 7690.067
$ws.Range("I2").Value = 13421
$ws.Range("J2").Value = 1140.4286
$ws.Range("K2").Value = 13421
$ws.Range("L2").Value = 1140.4286
$ws.Range("M2").Value = -13308
$ws.Range("N2").Value = -1366.4286
# Row 74
$ws.Range("H74").Value = 5647.12
$ws.Range("I74").Value = 875
$ws.Range("K74").Value = 875
$ws.Range("M74").Value = -1
# Row 77
$ws.Range("H77").Value = 5647.12
$ws.Range("I77").Value = 875
$ws.Range("K77").Value = 4375
$ws.Range("M77").Value = -7
# Row 110
$ws.Range("H110").Value = 1068.5217
$ws.Range("I110").Value = 977.6842
$ws.Range("K110").Value = 977.6842
$ws.Range("M110").Value = 1067.3158
# Row 116
$ws.Range("H116").Value = 7690.067
$ws.Range("I116").Value = 13421
$ws.Range("J116").Value = 1140.4286
$ws.Range("K116").Value = 13421
$ws.Range("L116").Value = 1140.4286
$ws.Range("M116").Value = -11127
$ws.Range("N116").Value = -5728.4286
# Row 122
$ws.Range("H122").Value = 13911.375
$ws.Range("I122").Value = 25772.25
$ws.Range("J122").Value = 2050.5
$ws.Range("K122").Value = 77316.75
$ws.Range("L122").Value = 6151.5
$ws.Range("M122").Value = -74866.75
$ws.Range("N122").Value = -11051.5
# Row 132
$ws.Range("H132").Value = 2471.0205
$ws.Range("I132").Value = 1873.7222
$ws.Range("K132").Value = 5621.1666
$ws.Range("M132").Value = -3091.1666
# Row 139
$ws.Range("H139").Value = 55857.5
$ws.Range("J139").Value = 55857.5
$ws.Range("L139").Value = 55857.5
$ws.Range("N139").Value = -66137.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 7690.067
$ws.Range("I3").Value = 13421
$ws.Range("J3").Value = 1140.4286
$ws.Range("K3").Value = 13421
$ws.Range("L3").Value = 1140.4286
$ws.Range("M3").Value = -13307
$ws.Range("N3").Value = -1368.4286
# Row 134
$ws.Range("H134").Value = 3267.9395
$ws.Range("I134").Value = 2156.739
$ws.Range("J134").Value = 5823.7
$ws.Range("K134").Value = 6470.217000000001
$ws.Range("L134").Value = 17471.1
$ws.Range("M134").Value = -3935.217000000001
$ws.Range("N134").Value = -22541.1

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 712.5714
$ws.Range("I16").Value = 720
$ws.Range("J16").Value = 694
$ws.Range("K16").Value = 720
$ws.Range("L16").Value = 694
$ws.Range("M16").Value = -433
$ws.Range("N16").Value = -1268
# Row 31
$ws.Range("H31").Value = 9964
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
# Row 34
$ws.Range("H34").Value = 9964
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
# Row 105
$ws.Range("H105").Value = 496.66666
$ws.Range("I105").Value = 556
$ws.Range("K105").Value = 556
$ws.Range("M105").Value = 1191
# Row 113
$ws.Range("H113").Value = 712.5714
$ws.Range("I113").Value = 720
$ws.Range("J113").Value = 694
$ws.Range("K113").Value = 720
$ws.Range("L113").Value = 694
$ws.Range("M113").Value = 1450
$ws.Range("N113").Value = -5034
# Row 122
$ws.Range("H122").Value = 1528
$ws.Range("I122").Value = 1037.3334
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 3112.0002
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -662.0001999999999
$ws.Range("N122").Value = -13900
# Row 132
$ws.Range("H132").Value = 3649.8333
$ws.Range("I132").Value = 2433.2
$ws.Range("J132").Value = 5170.625
$ws.Range("K132").Value = 7299.599999999999
$ws.Range("L132").Value = 15511.875
$ws.Range("M132").Value = -4769.599999999999
$ws.Range("N132").Value = -20571.875
# Row 134
$ws.Range("H134").Value = 4355.6816
$ws.Range("I134").Value = 2571.5
$ws.Range("J134").Value = 7478
$ws.Range("K134").Value = 7714.5
$ws.Range("L134").Value = 22434
$ws.Range("M134").Value = -5179.5
$ws.Range("N134").Value = -27504

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 49
$ws.Range("H49").Value = 1566
$ws.Range("J49").Value = 1566
$ws.Range("L49").Value = 4698
$ws.Range("N49").Value = -5010
# Row 70
$ws.Range("H70").Value = 2671.2856
$ws.Range("I70").Value = 566.3333
$ws.Range("J70").Value = 4250
$ws.Range("K70").Value = 1698.9999
$ws.Range("L70").Value = 12750
$ws.Range("M70").Value = -1383.9999
$ws.Range("N70").Value = -13380
# Row 73
$ws.Range("H73").Value = 2671.2856
$ws.Range("I73").Value = 566.3333
$ws.Range("J73").Value = 4250
$ws.Range("K73").Value = 1698.9999
$ws.Range("L73").Value = 12750
$ws.Range("M73").Value = -606.9999
$ws.Range("N73").Value = -14934
# Row 75
$ws.Range("H75").Value = 2863.2856
$ws.Range("J75").Value = 2983.4
$ws.Range("L75").Value = 8950.200000000001
$ws.Range("N75").Value = -10946.2
# Row 78
$ws.Range("H78").Value = 2863.2856
$ws.Range("J78").Value = 2983.4
$ws.Range("L78").Value = 26850.6
$ws.Range("N78").Value = -36834.60000000001
# Row 131
$ws.Range("H131").Value = 1254.39
$ws.Range("I131").Value = 235
$ws.Range("J131").Value = 1343.0326
$ws.Range("K131").Value = 705
$ws.Range("L131").Value = 4029.0978
$ws.Range("M131").Value = 4335
$ws.Range("N131").Value = -14109.0978

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1694.25
$ws.Range("I102").Value = 920.8
$ws.Range("J102").Value = 2983.3333
$ws.Range("K102").Value = 920.8
$ws.Range("L102").Value = 2983.3333
$ws.Range("M102").Value = 701.2
$ws.Range("N102").Value = -6227.3333
# Row 113
$ws.Range("H113").Value = 1200
$ws.Range("I113").Value = 1400
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 1400
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 770
$ws.Range("N113").Value = -5340
# Row 122
$ws.Range("H122").Value = 1391388.6
$ws.Range("I122").Value = 2224421.8
$ws.Range("K122").Value = 6673265.399999999
$ws.Range("M122").Value = -6670815.399999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 2436.2
$ws.Range("I46").Value = 2001
$ws.Range("J46").Value = 2545
$ws.Range("K46").Value = 2001
$ws.Range("L46").Value = 2545
$ws.Range("M46").Value = -1813
$ws.Range("N46").Value = -2921
# Row 122
$ws.Range("H122").Value = 3911.7646
# Row 132
$ws.Range("H132").Value = 4180
$ws.Range("I132").Value = 2715.4443
$ws.Range("J132").Value = 6576.5454
$ws.Range("K132").Value = 8146.3329
$ws.Range("L132").Value = 19729.6362
$ws.Range("M132").Value = -5616.3329
$ws.Range("N132").Value = -24789.6362
# Row 136
$ws.Range("H136").Value = 5444.375
$ws.Range("I136").Value = 2086.5715
$ws.Range("J136").Value = 6827
$ws.Range("K136").Value = 6259.7145
$ws.Range("L136").Value = 20481
$ws.Range("M136").Value = -3709.7145
$ws.Range("N136").Value = -25581

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 45380.566
$ws.Range("I122").Value = 101275.8
$ws.Range("J122").Value = 2384.2307
$ws.Range("K122").Value = 303827.4
$ws.Range("L122").Value = 7152.6921
$ws.Range("M122").Value = -301377.4
$ws.Range("N122").Value = -12052.6921
# Row 136
$ws.Range("H136").Value = 12860875
$ws.Range("I136").Value = 17597502
$ws.Range("K136").Value = 52792506
$ws.Range("M136").Value = -52789956

